# Updated cryptos list on Wed Mar 20 04:37:14 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make sure the Price column (D) keeps being stored as text, even for
# values that look numeric (e.g. "1.00", "10.60"), so trailing zeros and
# formatting are preserved exactly like the source data feed produces.
$ws.Range("D2:D51").NumberFormat = "@"

# --- Straightforward price / volume refreshes (rows whose coin stays put) ---
$ws.Range("D2").Value  = "61.562.63"
$ws.Range("E2").Value  = "  -6.32%  "

$ws.Range("D3").Value  = "3.133.51"
$ws.Range("E3").Value  = "  -8.39%  "

$ws.Range("E4").Value  = "  +0.37%  "

$ws.Range("D5").Value  = "508.18"
$ws.Range("E5").Value  = "  -5.11%  "

$ws.Range("D6").Value  = "167.08"
$ws.Range("E6").Value  = "  -11.36%  "

$ws.Range("D7").Value  = "0.578"
$ws.Range("E7").Value  = "  -5.73%  "

$ws.Range("E8").Value  = "  +0.28%  "

$ws.Range("D9").Value  = "3.142.66"
$ws.Range("E9").Value  = "  -8.00%  "

$ws.Range("D10").Value = "0.582"
$ws.Range("E10").Value = "  -8.46%  "

$ws.Range("D11").Value = "50.94"
$ws.Range("E11").Value = "  -15.22%  "

$ws.Range("D12").Value = "0.126"
$ws.Range("E12").Value = "  -7.41%  "

$ws.Range("D13").Value = "0.0000245"
$ws.Range("E13").Value = "  -5.58%  "

$ws.Range("D14").Value = "8.63"
$ws.Range("E14").Value = "  -8.28%  "

$ws.Range("D15").Value = "3.679.52"
$ws.Range("E15").Value = "  -6.79%  "

$ws.Range("D16").Value = "3.168.50"
$ws.Range("E16").Value = "  -6.93%  "

$ws.Range("D17").Value = "0.112"
$ws.Range("E17").Value = "  -8.85%  "

$ws.Range("D18").Value = "61.725.20"
$ws.Range("E18").Value = "  -5.63%  "

$ws.Range("D19").Value = "16.67"
$ws.Range("E19").Value = "  -5.98%  "

$ws.Range("D20").Value = "10.60"
$ws.Range("E20").Value = "  -6.63%  "

$ws.Range("D21").Value = "0.931"
$ws.Range("E21").Value = "  -5.87%  "

$ws.Range("D22").Value = "354.99"
$ws.Range("E22").Value = "  -6.16%  "

$ws.Range("D23").Value = "3.62"
$ws.Range("E23").Value = "  -4.80%  "

$ws.Range("D24").Value = "78.60"
$ws.Range("E24").Value = "  -5.23%  "

$ws.Range("D25").Value = "10.69"
$ws.Range("E25").Value = "  -3.16%  "

$ws.Range("D26").Value = "6.10"
$ws.Range("E26").Value = "  +3.82%  "

$ws.Range("D27").Value = "3.76"
$ws.Range("E27").Value = "  +0.86%  "

$ws.Range("D28").Value = "2.52"
$ws.Range("E28").Value = "  -6.52%  "

$ws.Range("D29").Value = "10.81"
$ws.Range("E29").Value = "  -9.08%  "

$ws.Range("D30").Value = "7.88"
$ws.Range("E30").Value = "  -9.18%  "

$ws.Range("D31").Value = "633.57"
$ws.Range("E31").Value = "  -9.37%  "

$ws.Range("D32").Value = "27.44"
$ws.Range("E32").Value = "  -8.70%  "

$ws.Range("D33").Value = "6.19"
$ws.Range("E33").Value = "  -9.93%  "

$ws.Range("D34").Value = "10.89"
$ws.Range("E34").Value = "  -4.35%  "

# --- Rows 35-37: the ranking reshuffled (Hedera dropped, Dai moved up) ---
# Old: 35=Hedera, 36=OKB, 37=Dai  ->  New: 35=Dai, 36=Hedera, 37=OKB
$ws.Range("B35").Value = "Dai"
$ws.Range("C35").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  -0.01%  "

$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "0.100"
$ws.Range("E36").Value = "  -6.75%  "

$ws.Range("B37").Value = "OKB"
$ws.Range("C37").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D37").Value = "56.08"
$ws.Range("E37").Value = "  -9.13%  "

# --- More straightforward refreshes ---
$ws.Range("D38").Value = "35.23"
$ws.Range("E38").Value = "  -5.07%  "

$ws.Range("D39").Value = "0.361"
$ws.Range("E39").Value = "  -7.11%  "

$ws.Range("E40").Value = "  +0.24%  "

$ws.Range("D41").Value = "0.0₃0670"
$ws.Range("E41").Value = "  +5.65%  "

$ws.Range("E42").Value = "  -7.39%  "

$ws.Range("D43").Value = "2.794.04"
$ws.Range("E43").Value = "  -3.92%  "

$ws.Range("D44").Value = "2.41"
$ws.Range("E44").Value = "  +0.39%  "

$ws.Range("E45").Value = "  -3.49%  "

# --- Rows 46-47: Stacks and VeChain swap places ---
# Old: 46=VeChain, 47=Stacks  ->  New: 46=Stacks, 47=VeChain
$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D46").Value = "2.72"
$ws.Range("E46").Value = "  +1.80%  "

$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").Value = "0.0374"
$ws.Range("E47").Value = "  -7.59%  "

# --- Final straightforward refreshes ---
$ws.Range("D48").Value = "2.48"
$ws.Range("E48").Value = "  -12.42%  "

$ws.Range("D49").Value = "2.87"
$ws.Range("E49").Value = "  -0.57%  "

$ws.Range("D50").Value = "132.52"
$ws.Range("E50").Value = "  -4.05%  "

$ws.Range("E51").Value = "  -6.33%  "
